$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some look numeric, e.g. '214.65',
# others use '.' as a thousands separator, e.g. '27.009.32', and are
# never valid numbers). The source workbook always stores these as
# plain text (inline/shared strings), never as numbers. Force the
# Price column to Text *before* writing so Excel doesn't reinterpret
# numeric-looking values (losing exact digits to float rounding),
# then restore the default (Normal) style so no stray per-cell
# formatting is left behind once the values are set.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "27.009.32"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.620.43"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "214.65"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "0.0628"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.847.69"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "1.629.48"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "64.76"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("D17").Value = "26.985.83"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "0.0₃0746"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "213.79"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "2.38"
$ws.Range("E23").Value = "  -5.31%  "
$ws.Range("D24").Value = "9.04"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "148.14"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "0.117"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").Value = "0.752"
$ws.Range("E33").Value = "  +36.60%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "1.345.55"
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "0.0177"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "0.846"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "0.800"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "65.01"
$ws.Range("E43").Value = "  +5.21%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "1.758.48"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("D51").Value = "0.0515"
$ws.Range("E51").Value = "  +0.39%  "

# Rows 46 and 47: Quant and WEMIXToken swap list positions
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "0.881"
$ws.Range("E46").Value = "  +31.57%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "89.76"
$ws.Range("E47").Value = "  -2.19%  "

# Drop the temporary Text number format again so styling matches
# the original (untouched) workbook.
$priceCol.Style = "Normal"
